# Updated symbol list on Mon Jan  9 19:24:42 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for the crypto rows whose
# market data moved since the previous run. Values are written with a
# leading quote so Excel keeps them as literal text (matching the sheet's
# existing inline-string formatting like "7.000" / "4.72%") instead of
# coercing them into numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'278.45"
$ws.Range("E2").Value = "'4.78%"

$ws.Range("D3").Value = "'26.88"
$ws.Range("E3").Value = "'0.69%"

$ws.Range("D4").Value = "'4.920"
$ws.Range("E4").Value = "'4.48%"

$ws.Range("D5").Value = "'0.06393"
$ws.Range("E5").Value = "'5.16%"

$ws.Range("D6").Value = "'7.002"
$ws.Range("E6").Value = "'3.86%"

$ws.Range("D7").Value = "'3.356"
$ws.Range("E7").Value = "'5.89%"

$ws.Range("D8").Value = "'0.8885"
$ws.Range("E8").Value = "'4.55%"

$ws.Range("D9").Value = "'1.146"
$ws.Range("E9").Value = "'26.61%"

$ws.Range("D10").Value = "'0.1493"
$ws.Range("E10").Value = "'5.92%"

$ws.Range("D11").Value = "'0.05256"
$ws.Range("E11").Value = "'8.02%"

$ws.Range("D12").Value = "'0.07331"
$ws.Range("E12").Value = "'3.42%"

$ws.Range("D13").Value = "'0.03118"
$ws.Range("E13").Value = "'-1.91%"

$ws.Range("D14").Value = "'0.09064"
$ws.Range("E14").Value = "'0.46%"

$ws.Range("D15").Value = "'0.001568"
$ws.Range("E15").Value = "'2.25%"

$ws.Range("D16").Value = "'0.0006342"
$ws.Range("E16").Value = "'4.72%"

$ws.Range("D17").Value = "'0.006038"
$ws.Range("E17").Value = "'0.62%"

$ws.Range("D18").Value = "'3.491"
$ws.Range("E18").Value = "'1.00%"

$ws.Range("D19").Value = "'2.281"
$ws.Range("E19").Value = "'0.16%"

$ws.Range("D20").Value = "'0.3157"
$ws.Range("E20").Value = "'3.33%"

$ws.Range("E21").Value = "'2.52%"

$ws.Range("D22").Value = "'3.927"
$ws.Range("E22").Value = "'-3.68%"

$ws.Range("D23").Value = "'0.04364"
$ws.Range("E23").Value = "'2.80%"

$ws.Range("D24").Value = "'0.001181"
$ws.Range("E24").Value = "'-0.30%"

$ws.Range("D25").Value = "'0.003679"
$ws.Range("E25").Value = "'-11.05%"

$ws.Range("E26").Value = "'0.09%"

$ws.Range("D27").Value = "'0.0001700"
$ws.Range("E27").Value = "'1.12%"

$ws.Range("D40").Value = "'0.04077"
$ws.Range("E40").Value = "'4.22%"

$ws.Range("D41").Value = "'0.006648"
$ws.Range("E41").Value = "'59.54%"

$ws.Range("E42").Value = "'5.42%"

$ws.Range("D43").Value = "'0.002363"
$ws.Range("E43").Value = "'11.96%"

$ws.Range("E44").Value = "'2.23%"

$ws.Range("D45").Value = "'0.00005269"
$ws.Range("E45").Value = "'3.18%"

$ws.Range("E46").Value = "'-0.02%"

$ws.Range("E47").Value = "'1,406.78%"

$ws.Range("D48").Value = "'0.02123"
$ws.Range("E48").Value = "'-13.28%"

$ws.Range("E49").Value = "'-0.02%"

$ws.Range("E50").Value = "'-0.09%"
